$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 574, shifting existing rows 574-634 down to 575-635
$ws.Rows.Item(574).Insert()

$ws.Cells.Item(574, 1).Value = 5
$ws.Cells.Item(574, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(574, 3).Value = "Maule"
$ws.Cells.Item(574, 4).Value = 45194
$ws.Cells.Item(574, 5).Value = 7
$ws.Cells.Item(574, 6).Value = 100112032
$ws.Cells.Item(574, 7).Value = "Zapallo italiano"
$ws.Cells.Item(574, 8).Value = "Sin especificar"
$ws.Cells.Item(574, 9).Value = "Primera"
$ws.Cells.Item(574, 10).Value = 300
$ws.Cells.Item(574, 11).Value = 12000
$ws.Cells.Item(574, 12).Value = 12000
$ws.Cells.Item(574, 13).Value = 12000
$ws.Cells.Item(574, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(574, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(574, 16).Value = 240
$ws.Cells.Item(574, 17).Value = 50
$ws.Cells.Item(574, 18).Value = "Hortaliza"
